# Actualizacion Estado de Cuenta - NIT 9008623774
# Elimina los periodos de mora anteriores y los reemplaza por los
# nuevos periodos de mora (se reconstruye la base de datos de la hoja).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

function Get-DescendingPeriods {
    param($StartYear, $StartMonth, $Count)
    $result = @()
    $yy = $StartYear
    $mm = $StartMonth
    for ($i = 0; $i -lt $Count; $i++) {
        $result += ("{0:D2}{1:D2}" -f $yy, $mm)
        $mm = $mm - 1
        if ($mm -eq 0) {
            $mm = 12
            $yy = $yy - 1
        }
    }
    return $result
}

# Trabajadores con mora vigente (NIT del empleador: 9008623774)
$idYulis   = "1143356636"
$nomYulis  = "YULIS PATRICIA HUERTAS RODRIGUEZ"
$salYulis  = 782000

$idYajaira  = "1047469042"
$nomYajaira = "YAJAIRA ISABEL GAMARRA ESPITIA"
$salYajaira = 781242

$idJorge   = "1047365583"
$nomJorge  = "JORGE ARMANDO BARRIOS CASTRO"
$salJorge  = 781242

$idHainner  = "1034278383"
$nomHainner = "HAINNER JOSE CHACON ROMERO"
$salHainner = 782000

# Periodos adeudados por cada trabajador (de mas reciente a mas antiguo)
$periodosJorge   = Get-DescendingPeriods 21 5 37
$periodosHainner = Get-DescendingPeriods 21 5 39

# Construye las filas finales (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
$rows = @()
$rows += ,@("CC", $idYulis, $nomYulis, "1803", 31280, $salYulis)
$rows += ,@("CC", $idYajaira, $nomYajaira, "1805", 10416, $salYajaira)

for ($i = 0; $i -lt $periodosJorge.Count; $i++) {
    $periodo = $periodosJorge[$i]
    $valor = 31249
    if ($periodo -eq "2105") { $valor = 24999 }
    if ($periodo -eq "1805") { $valor = 7291 }
    $rows += ,@("CC", $idJorge, $nomJorge, $periodo, $valor, $salJorge)
}

for ($i = 0; $i -lt $periodosHainner.Count; $i++) {
    $periodo = $periodosHainner[$i]
    $valor = 31280
    if ($periodo -eq "2105") { $valor = 25024 }
    if ($periodo -eq "1803") { $valor = 16683 }
    $rows += ,@("CC", $idHainner, $nomHainner, $periodo, $valor, $salHainner)
}

# Vuelca las filas nuevas sobre el rango de datos existente (filas 16 a 93)
$firstRow = 16
$numRows = $rows.Count
$arr = New-Object 'object[,]' $numRows,6
for ($r = 0; $r -lt $numRows; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt 6; $c++) {
        $arr[$r, $c] = $row[$c]
    }
}

$lastRow = $firstRow + $numRows - 1
$targetRange = $ws.Range("B$firstRow" + ":G$lastRow")
$targetRange.Value = $arr

# Reajusta el ancho de columnas al nuevo contenido
$ws.Columns.Item("B:J").AutoFit() | Out-Null

$wb.Save()
